# Adds the two new matches (rows 102 and 103) that were scraped for the
# ecuador / liga-pro / 2023 sheet, matching the "Atualizado por script em
# 12-11-2023 20:45" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New row data, in column order A..V (Indice is filled in separately as
# it is simply row-number - 1).
# ---------------------------------------------------------------------
$newRows = @(
    @{
        Row = 102
        Indice = 101
        Pais = "ecuador"
        Torneio = "liga-pro"
        Temporada = "2023"
        DataPartida = 45242.77083333334
        Home = "EL Nacional"
        HomeGols = 4
        Away = "Gualaceo"
        AwayGols = 2
        HomeOpenOdds = 1.45
        HomeOpenData = "06/11/2023 00:12"
        HomeCloseOdds = 1.48
        HomeCloseData = "12/11/2023 18:29"
        DrawOpenOdds = 4.64
        DrawOpenData = "06/11/2023 00:12"
        DrawCloseOdds = 4.59
        DrawCloseData = "12/11/2023 18:29"
        AwayOpenOdds = 6.35
        AwayOpenData = "06/11/2023 00:12"
        AwayCloseOdds = 6.31
        AwayCloseData = "12/11/2023 18:29"
        Url = "https://www.betexplorer.com/football/ecuador/liga-pro/el-nacional-gualaceo/8UFjCyl2/"
    },
    @{
        Row = 103
        Indice = 102
        Pais = "ecuador"
        Torneio = "liga-pro"
        Temporada = "2023"
        DataPartida = 45242.875
        Home = "Libertad"
        HomeGols = 0
        Away = "LDU Quito"
        AwayGols = 2
        HomeOpenOdds = 4.1
        HomeOpenData = "09/11/2023 01:12"
        HomeCloseOdds = 5.18
        HomeCloseData = "12/11/2023 20:58"
        DrawOpenOdds = 3.64
        DrawOpenData = "09/11/2023 01:12"
        DrawCloseOdds = 3.29
        DrawCloseData = "12/11/2023 20:58"
        AwayOpenOdds = 1.79
        AwayOpenData = "09/11/2023 01:12"
        AwayCloseOdds = 1.82
        AwayCloseData = "12/11/2023 20:58"
        Url = "https://www.betexplorer.com/football/ecuador/liga-pro/libertad-ldu-quito/Ykd2YR44/"
    }
)

# Template cells from the last existing data row (101): carries the
# correct number formats / bold+border style for the "Indice" column (A)
# and the date-time style for "data_partida" (E).
$indiceTemplate = $ws.Range("A101")
$dataTemplate = $ws.Range("E101")

foreach ($r in $newRows) {
    $row = $r.Row

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $r.Indice
    $cellA.NumberFormat = $indiceTemplate.NumberFormat
    $cellA.Font.Bold = $true
    $cellA.HorizontalAlignment = $indiceTemplate.HorizontalAlignment
    $cellA.VerticalAlignment = $indiceTemplate.VerticalAlignment
    $cellA.Borders.LineStyle = $indiceTemplate.Borders.LineStyle

    $ws.Cells.Item($row, 2).Value = $r.Pais
    $ws.Cells.Item($row, 3).Value = $r.Torneio
    # "temporada" is a digit-only string ("2023"); prefix with an
    # apostrophe so the engine keeps it text instead of coercing it to a
    # number (Excel strips the leading quote from the stored value).
    $ws.Cells.Item($row, 4).Value = "'" + $r.Temporada

    $cellE = $ws.Cells.Item($row, 5)
    $cellE.Value = $r.DataPartida
    $cellE.NumberFormat = $dataTemplate.NumberFormat

    $ws.Cells.Item($row, 6).Value = $r.Home
    $ws.Cells.Item($row, 7).Value = $r.HomeGols
    $ws.Cells.Item($row, 8).Value = $r.Away
    $ws.Cells.Item($row, 9).Value = $r.AwayGols

    $ws.Cells.Item($row, 10).Value = $r.HomeOpenOdds
    $ws.Cells.Item($row, 11).Value = $r.HomeOpenData
    $ws.Cells.Item($row, 12).Value = $r.HomeCloseOdds
    $ws.Cells.Item($row, 13).Value = $r.HomeCloseData

    $ws.Cells.Item($row, 14).Value = $r.DrawOpenOdds
    $ws.Cells.Item($row, 15).Value = $r.DrawOpenData
    $ws.Cells.Item($row, 16).Value = $r.DrawCloseOdds
    $ws.Cells.Item($row, 17).Value = $r.DrawCloseData

    $ws.Cells.Item($row, 18).Value = $r.AwayOpenOdds
    $ws.Cells.Item($row, 19).Value = $r.AwayOpenData
    $ws.Cells.Item($row, 20).Value = $r.AwayCloseOdds
    $ws.Cells.Item($row, 21).Value = $r.AwayCloseData

    $ws.Cells.Item($row, 22).Value = $r.Url
}

Write-Host "Added rows 102:103 to Sheet1"
